$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linking option C test fixed: swap in the correct linked account
# (sa021@mailinator.com -> sa019@mailinator.com) for row 6.
$ws.Range("B6").Value = "sa019@mailinator.com"

# Record the (Mailing Address / Residential Address) for that same row,
# same text re-used in both columns M and N.
$ws.Range("M6").Value = "UNIT 35, 146-152 PARRAMATTA RD, HOMEBUSH NSW 2140"
$ws.Range("N6").Value = "UNIT 35, 146-152 PARRAMATTA RD, HOMEBUSH NSW 2140"

# Scroll the view so column C is the left-most visible column and select N6.
$excel.ActiveWindow.ScrollColumn = 3
[void]$ws.Range("N6").Select()

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
